$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new ProductWarranty column (Z) with its header and data value
$ws.Range("Z1").Value = "ProductWarranty"
$ws.Range("Z2").Value = "a2K1g000000CgAT"

# Give column Z the same look as the other data columns (width ~28 chars)
$ws.Columns.Item(26).ColumnWidth = 27.17

# Move the view so column Z (around column T) is visible, and select the new bottom-right area
$excel.ActiveWindow.ScrollColumn = 20
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("X15").Select()
